# motive analysis phase 1
# Add a new "murders per capita (population-normalized)" column M:
#   m_p_pop = (sum of yearly murder counts 2015-2021) / (region pop * 7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in M1
$ws.Range("M1").Value = "m_p_pop"

# M2 gets its own (non-shared) formula
$ws.Range("M2").Formula = "=(D2+E2+F2+G2+H2+I2+J2)/(L2*7)"

# M3:M37 share one relative formula (Excel will store this as a shared
# formula group when the same formula text is assigned to a multi-cell
# range in one shot)
$ws.Range("M3:M37").Formula = "=(D3+E3+F3+G3+H3+I3+J3)/(L3*7)"

# Restore the selection to L16 as left by the author after this edit
$ws.Range("L16").Select() | Out-Null
